{"js": "// Fill in the \"Expected Effort\" / \"Actual Effort\" columns of the Project\n// Schedule table with the estimated / available hours, and clean up the\n// stray \"_GoBack\" bookmark left over from the previous save.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nasync function setCellText(rowIndex, colIndex, text) {\n  const cell = table.getCell(rowIndex, colIndex);\n  cell.body.clear();\n  cell.body.insertText(text, \"Replace\");\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n  cell.body.paragraphs.items[0].alignment = \"Centered\";\n  await context.sync();\n}\n\nasync function centerCell(rowIndex, colIndex) {\n  const cell = table.getCell(rowIndex, colIndex);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n  cell.body.paragraphs.items[0].alignment = \"Centered\";\n  await context.sync();\n}\n\n// Row 1 = \"Timeline\"\nawait setCellText(1, 1, \"2\");\nawait setCellText(1, 2, \"1\");\n\n// Row 2 = \"SIQ\"\nawait setCellText(2, 1, \"1\");\nawait setCellText(2, 2, \"1\");\n\n// Row 3 = \"PMP\"\nawait setCellText(3, 1, \"10\");\nawait setCellText(3, 2, \"10\");\n\n// Row 4 = \"SRS\"\nawait setCellText(4, 1, \"12\");\nawait setCellText(4, 2, \"10\");\n\n// Row 5 (trailing blank row) - only the first (label) cell gets centered,\n// no text is added.\nawait centerCell(5, 0);\n\n// Remove the leftover \"_GoBack\" bookmark (Word drops this automatically on\n// a fresh save once the cursor position it tracked is no longer relevant).\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // Bookmark may already be absent; nothing else to do.\n}\n", "ps1": "# Fill in the \"Expected Effort\" / \"Actual Effort\" columns of the Project\n# Schedule table with the estimated / available hours, and clean up the\n# stray \"_GoBack\" bookmark left over from the previous save.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Set-CellValue($row, $col, $text) {\n    $cell = $tbl.Cell($row, $col)\n    $cell.Range.Text = $text\n    $cell.Range.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter\n}\n\n# Row 2 = \"Timeline\"\nSet-CellValue 2 2 \"2\"\nSet-CellValue 2 3 \"1\"\n\n# Row 3 = \"SIQ\"\nSet-CellValue 3 2 \"1\"\nSet-CellValue 3 3 \"1\"\n\n# Row 4 = \"PMP\"\nSet-CellValue 4 2 \"10\"\nSet-CellValue 4 3 \"10\"\n\n# Row 5 = \"SRS\"\nSet-CellValue 5 2 \"12\"\nSet-CellValue 5 3 \"10\"\n\n# Row 6 (trailing blank row) - only the first (label) cell gets centered,\n# no text is added.\n$lastCell = $tbl.Cell(6, 1)\n$lastCell.Range.ParagraphFormat.Alignment = 1\n\n# Remove the leftover \"_GoBack\" bookmark (Word drops this automatically on\n# a fresh save once the cursor position it tracked is no longer relevant).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
